# Updated cryptos list on Thu May 16 22:41:22 UTC 2024 with GitHub Actions
#
# Refreshes the "Price" (column D) and "Volume(1h)" (column E) columns of the
# crypto table with newly scraped values. All cells in these columns are
# plain text (e.g. "65.324.74", "  -1.30%  "), not real numbers, so every
# write must preserve the original text representation exactly (including
# values that look numeric, like "6.72" or "0.300") and must not leave the
# cell's style/number-format altered.
#
# For column D, values such as "6.72" or "0.0000101" would otherwise be
# auto-coerced by Excel into floating-point numbers (losing trailing zeros /
# exact text), so we temporarily force a text number-format before writing,
# then clear formatting again so the cell's style index is left untouched
# (matches the original workbook, which has no explicit style on these
# cells). Column E values always contain spaces/percent signs, so they are
# never reinterpreted as numbers and can be set directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $cell = $ws.Range($range)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

function Set-PlainValue($range, $value) {
    $ws.Range($range).Value = $value
}

# row -> (new Price, new Volume(1h)); Price omitted where unchanged
$rowUpdates = @(
    @{ Row = 2;  Price = "65.324.74";  Volume = "  -1.30%  " }
    @{ Row = 3;  Price = "2.933.05";   Volume = "  -2.82%  " }
    @{ Row = 4;  Price = $null;        Volume = "  -0.13%  " }
    @{ Row = 5;  Price = "568.16";     Volume = "  -2.98%  " }
    @{ Row = 6;  Price = "157.97";     Volume = "  +1.09%  " }
    @{ Row = 7;  Price = $null;        Volume = "  +0.05%  " }
    @{ Row = 8;  Price = $null;        Volume = "  -0.64%  " }
    @{ Row = 9;  Price = "2.927.76";   Volume = "  -2.83%  " }
    @{ Row = 10; Price = "6.72";       Volume = "  -3.53%  " }
    @{ Row = 11; Price = $null;        Volume = "  -3.88%  " }
    @{ Row = 12; Price = "0.458";      Volume = "  +1.32%  " }
    @{ Row = 13; Price = "0.0000244";  Volume = "  -3.14%  " }
    @{ Row = 14; Price = "34.25";      Volume = "  -1.16%  " }
    @{ Row = 15; Price = $null;        Volume = "  -0.81%  " }
    @{ Row = 16; Price = "65.274.05";  Volume = "  -1.39%  " }
    @{ Row = 17; Price = "3.417.45";   Volume = "  -2.90%  " }
    @{ Row = 18; Price = $null;        Volume = "  +0.14%  " }
    @{ Row = 19; Price = "2.930.73";   Volume = "  -2.94%  " }
    @{ Row = 20; Price = "15.63";      Volume = "  +12.64%  " }
    @{ Row = 21; Price = "443.36";     Volume = "  -4.32%  " }
    @{ Row = 22; Price = $null;        Volume = "  +0.74%  " }
    @{ Row = 23; Price = "7.26";       Volume = "  -1.48%  " }
    @{ Row = 24; Price = "82.11";      Volume = "  +0.05%  " }
    @{ Row = 25; Price = "2.23";       Volume = "  -1.43%  " }
    @{ Row = 26; Price = "12.09";      Volume = "  -3.40%  " }
    @{ Row = 27; Price = "10.06";      Volume = "  -6.06%  " }
    @{ Row = 28; Price = $null;        Volume = "  +0.07%  " }
    @{ Row = 29; Price = $null;        Volume = "  -0.49%  " }
    @{ Row = 30; Price = $null;        Volume = "  -0.04%  " }
    @{ Row = 31; Price = $null;        Volume = "  -1.64%  " }
    @{ Row = 32; Price = "0.0000101";  Volume = "  -4.39%  " }
    @{ Row = 33; Price = $null;        Volume = "  +0.00%  " }
    @{ Row = 34; Price = $null;        Volume = "  -0.32%  " }
    @{ Row = 35; Price = $null;        Volume = "  +0.04%  " }
    @{ Row = 36; Price = "0.970";      Volume = "  -3.02%  " }
    @{ Row = 37; Price = "5.73";       Volume = "  -1.51%  " }
    @{ Row = 38; Price = "49.61";      Volume = "  +0.17%  " }
    @{ Row = 39; Price = "44.72";      Volume = "  +2.00%  " }
    @{ Row = 40; Price = $null;        Volume = "  -8.94%  " }
    @{ Row = 41; Price = "0.300";      Volume = "  -1.08%  " }
    @{ Row = 42; Price = $null;        Volume = "  -2.42%  " }
    @{ Row = 43; Price = $null;        Volume = "  -7.97%  " }
    @{ Row = 44; Price = $null;        Volume = "  +0.39%  " }
    @{ Row = 45; Price = "381.42";     Volume = "  -3.33%  " }
    @{ Row = 46; Price = "0.0351";     Volume = "  -0.99%  " }
    @{ Row = 47; Price = "2.697.12";   Volume = "  -3.63%  " }
    @{ Row = 48; Price = "133.45";     Volume = "  -0.36%  " }
    @{ Row = 50; Price = $null;        Volume = "  +4.32%  " }
)

foreach ($update in $rowUpdates) {
    if ($null -ne $update.Price) {
        Set-TextValue "D$($update.Row)" $update.Price
    }
    Set-PlainValue "E$($update.Row)" $update.Volume
}
